# Edit script: "6 month coverage from 2026 in scenario 3a"
# Expands the yearly coverage timeline (2018-2040) into a half-yearly
# timeline (2018, 2018.5, 2019, ... 2040) on both worksheets, and moves
# the scaled-up coverage (New Product A / the 0.8 & 0.5 age-band rows)
# onto 6-month steps starting in 2026.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Platform Coverage"
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Platform Coverage")

# Row 1: year headers H1:AZ1, half-yearly steps 2018 -> 2040
$ws1.Cells.Item(1, 8).Value = 2018
$ws1.Cells.Item(1, 9).Value = 2018.5
$ws1.Cells.Item(1, 10).Value = 2019
$ws1.Cells.Item(1, 11).Value = 2019.5
$ws1.Cells.Item(1, 12).Value = 2020
$ws1.Cells.Item(1, 13).Value = 2020.5
$ws1.Cells.Item(1, 14).Value = 2021
$ws1.Cells.Item(1, 15).Value = 2021.5
$ws1.Cells.Item(1, 16).Value = 2022
$ws1.Cells.Item(1, 17).Value = 2022.5
$ws1.Cells.Item(1, 18).Value = 2023
$ws1.Cells.Item(1, 19).Value = 2023.5
$ws1.Cells.Item(1, 20).Value = 2024
$ws1.Cells.Item(1, 21).Value = 2024.5
$ws1.Cells.Item(1, 22).Value = 2025
$ws1.Cells.Item(1, 23).Value = 2025.5
$ws1.Cells.Item(1, 24).Value = 2026
$ws1.Cells.Item(1, 25).Value = 2026.5
$ws1.Cells.Item(1, 26).Value = 2027
$ws1.Cells.Item(1, 27).Value = 2027.5
$ws1.Cells.Item(1, 28).Value = 2028
$ws1.Cells.Item(1, 29).Value = 2028.5
$ws1.Cells.Item(1, 30).Value = 2029
$ws1.Cells.Item(1, 31).Value = 2029.5
$ws1.Cells.Item(1, 32).Value = 2030
$ws1.Cells.Item(1, 33).Value = 2030.5
$ws1.Cells.Item(1, 34).Value = 2031
$ws1.Cells.Item(1, 35).Value = 2031.5
$ws1.Cells.Item(1, 36).Value = 2032
$ws1.Cells.Item(1, 37).Value = 2032.5
$ws1.Cells.Item(1, 38).Value = 2033
$ws1.Cells.Item(1, 39).Value = 2033.5
$ws1.Cells.Item(1, 40).Value = 2034
$ws1.Cells.Item(1, 41).Value = 2034.5
$ws1.Cells.Item(1, 42).Value = 2035
$ws1.Cells.Item(1, 43).Value = 2035.5
$ws1.Cells.Item(1, 44).Value = 2036
$ws1.Cells.Item(1, 45).Value = 2036.5
$ws1.Cells.Item(1, 46).Value = 2037
$ws1.Cells.Item(1, 47).Value = 2037.5
$ws1.Cells.Item(1, 48).Value = 2038
$ws1.Cells.Item(1, 49).Value = 2038.5
$ws1.Cells.Item(1, 50).Value = 2039
$ws1.Cells.Item(1, 51).Value = 2039.5
$ws1.Cells.Item(1, 52).Value = 2040

# Row 2 (0.6 coverage, age 5-15): add half-yearly continuation 2022-2025
$ws1.Range("P2").Value = 0.6
$ws1.Range("R2").Value = 0.6
$ws1.Range("T2").Value = 0.6
$ws1.Range("V2").Value = 0.6

# Row 3 (0.8 coverage, age 2-15): drop old biennial cells, cover every
# half-year from 2026 (X) through 2040 (AZ)
$ws1.Range("P3").ClearContents()
$ws1.Range("R3").ClearContents()
$ws1.Range("T3").ClearContents()
$ws1.Range("V3").ClearContents()
$ws1.Range("X3").Value = 0.8
$ws1.Range("Y3").Value = 0.8
$ws1.Range("Z3").Value = 0.8
$ws1.Range("AA3").Value = 0.8
$ws1.Range("AB3").Value = 0.8
$ws1.Range("AC3").Value = 0.8
$ws1.Range("AD3").Value = 0.8
$ws1.Range("AE3").Value = 0.8
$ws1.Range("AF3").Value = 0.8
$ws1.Range("AG3").Value = 0.8
$ws1.Range("AH3").Value = 0.8
$ws1.Range("AI3").Value = 0.8
$ws1.Range("AJ3").Value = 0.8
$ws1.Range("AK3").Value = 0.8
$ws1.Range("AL3").Value = 0.8
$ws1.Range("AM3").Value = 0.8
$ws1.Range("AN3").Value = 0.8
$ws1.Range("AO3").Value = 0.8
$ws1.Range("AP3").Value = 0.8
$ws1.Range("AQ3").Value = 0.8
$ws1.Range("AR3").Value = 0.8
$ws1.Range("AS3").Value = 0.8
$ws1.Range("AT3").Value = 0.8
$ws1.Range("AU3").Value = 0.8
$ws1.Range("AV3").Value = 0.8
$ws1.Range("AW3").Value = 0.8
$ws1.Range("AX3").Value = 0.8
$ws1.Range("AY3").Value = 0.8
$ws1.Range("AZ3").Value = 0.8

# Row 4 (0.5 coverage, age 15-50): drop old biennial cells, cover every
# half-year from 2026 (X) through 2040 (AZ)
$ws1.Range("P4").ClearContents()
$ws1.Range("R4").ClearContents()
$ws1.Range("T4").ClearContents()
$ws1.Range("V4").ClearContents()
$ws1.Range("X4").Value = 0.5
$ws1.Range("Y4").Value = 0.5
$ws1.Range("Z4").Value = 0.5
$ws1.Range("AA4").Value = 0.5
$ws1.Range("AB4").Value = 0.5
$ws1.Range("AC4").Value = 0.5
$ws1.Range("AD4").Value = 0.5
$ws1.Range("AE4").Value = 0.5
$ws1.Range("AF4").Value = 0.5
$ws1.Range("AG4").Value = 0.5
$ws1.Range("AH4").Value = 0.5
$ws1.Range("AI4").Value = 0.5
$ws1.Range("AJ4").Value = 0.5
$ws1.Range("AK4").Value = 0.5
$ws1.Range("AL4").Value = 0.5
$ws1.Range("AM4").Value = 0.5
$ws1.Range("AN4").Value = 0.5
$ws1.Range("AO4").Value = 0.5
$ws1.Range("AP4").Value = 0.5
$ws1.Range("AQ4").Value = 0.5
$ws1.Range("AR4").Value = 0.5
$ws1.Range("AS4").Value = 0.5
$ws1.Range("AT4").Value = 0.5
$ws1.Range("AU4").Value = 0.5
$ws1.Range("AV4").Value = 0.5
$ws1.Range("AW4").Value = 0.5
$ws1.Range("AX4").Value = 0.5
$ws1.Range("AY4").Value = 0.5
$ws1.Range("AZ4").Value = 0.5

# Row 5 (0.5 coverage, age 50-65): drop old biennial cells, cover every
# half-year from 2026 (X) through 2040 (AZ)
$ws1.Range("P5").ClearContents()
$ws1.Range("R5").ClearContents()
$ws1.Range("T5").ClearContents()
$ws1.Range("V5").ClearContents()
$ws1.Range("X5").Value = 0.5
$ws1.Range("Y5").Value = 0.5
$ws1.Range("Z5").Value = 0.5
$ws1.Range("AA5").Value = 0.5
$ws1.Range("AB5").Value = 0.5
$ws1.Range("AC5").Value = 0.5
$ws1.Range("AD5").Value = 0.5
$ws1.Range("AE5").Value = 0.5
$ws1.Range("AF5").Value = 0.5
$ws1.Range("AG5").Value = 0.5
$ws1.Range("AH5").Value = 0.5
$ws1.Range("AI5").Value = 0.5
$ws1.Range("AJ5").Value = 0.5
$ws1.Range("AK5").Value = 0.5
$ws1.Range("AL5").Value = 0.5
$ws1.Range("AM5").Value = 0.5
$ws1.Range("AN5").Value = 0.5
$ws1.Range("AO5").Value = 0.5
$ws1.Range("AP5").Value = 0.5
$ws1.Range("AQ5").Value = 0.5
$ws1.Range("AR5").Value = 0.5
$ws1.Range("AS5").Value = 0.5
$ws1.Range("AT5").Value = 0.5
$ws1.Range("AU5").Value = 0.5
$ws1.Range("AV5").Value = 0.5
$ws1.Range("AW5").Value = 0.5
$ws1.Range("AX5").Value = 0.5
$ws1.Range("AY5").Value = 0.5
$ws1.Range("AZ5").Value = 0.5

# ---------------------------------------------------------------
# Sheet "MarketShare"
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("MarketShare")

# Row 1: year headers D1:AV1, half-yearly steps 2018 -> 2040
$ws2.Cells.Item(1, 4).Value = 2018
$ws2.Cells.Item(1, 5).Value = 2018.5
$ws2.Cells.Item(1, 6).Value = 2019
$ws2.Cells.Item(1, 7).Value = 2019.5
$ws2.Cells.Item(1, 8).Value = 2020
$ws2.Cells.Item(1, 9).Value = 2020.5
$ws2.Cells.Item(1, 10).Value = 2021
$ws2.Cells.Item(1, 11).Value = 2021.5
$ws2.Cells.Item(1, 12).Value = 2022
$ws2.Cells.Item(1, 13).Value = 2022.5
$ws2.Cells.Item(1, 14).Value = 2023
$ws2.Cells.Item(1, 15).Value = 2023.5
$ws2.Cells.Item(1, 16).Value = 2024
$ws2.Cells.Item(1, 17).Value = 2024.5
$ws2.Cells.Item(1, 18).Value = 2025
$ws2.Cells.Item(1, 19).Value = 2025.5
$ws2.Cells.Item(1, 20).Value = 2026
$ws2.Cells.Item(1, 21).Value = 2026.5
$ws2.Cells.Item(1, 22).Value = 2027
$ws2.Cells.Item(1, 23).Value = 2027.5
$ws2.Cells.Item(1, 24).Value = 2028
$ws2.Cells.Item(1, 25).Value = 2028.5
$ws2.Cells.Item(1, 26).Value = 2029
$ws2.Cells.Item(1, 27).Value = 2029.5
$ws2.Cells.Item(1, 28).Value = 2030
$ws2.Cells.Item(1, 29).Value = 2030.5
$ws2.Cells.Item(1, 30).Value = 2031
$ws2.Cells.Item(1, 31).Value = 2031.5
$ws2.Cells.Item(1, 32).Value = 2032
$ws2.Cells.Item(1, 33).Value = 2032.5
$ws2.Cells.Item(1, 34).Value = 2033
$ws2.Cells.Item(1, 35).Value = 2033.5
$ws2.Cells.Item(1, 36).Value = 2034
$ws2.Cells.Item(1, 37).Value = 2034.5
$ws2.Cells.Item(1, 38).Value = 2035
$ws2.Cells.Item(1, 39).Value = 2035.5
$ws2.Cells.Item(1, 40).Value = 2036
$ws2.Cells.Item(1, 41).Value = 2036.5
$ws2.Cells.Item(1, 42).Value = 2037
$ws2.Cells.Item(1, 43).Value = 2037.5
$ws2.Cells.Item(1, 44).Value = 2038
$ws2.Cells.Item(1, 45).Value = 2038.5
$ws2.Cells.Item(1, 46).Value = 2039
$ws2.Cells.Item(1, 47).Value = 2039.5
$ws2.Cells.Item(1, 48).Value = 2040

# Row 2 (New Product A): drop old biennial cells (2026-2033), cover every
# half-year from 2026 (T) through 2040 (AV)
$ws2.Range("L2").ClearContents()
$ws2.Range("M2").ClearContents()
$ws2.Range("N2").ClearContents()
$ws2.Range("O2").ClearContents()
$ws2.Range("P2").ClearContents()
$ws2.Range("Q2").ClearContents()
$ws2.Range("R2").ClearContents()
$ws2.Range("S2").ClearContents()
$ws2.Range("T2").Value = 1
$ws2.Range("U2").Value = 1
$ws2.Range("V2").Value = 1
$ws2.Range("W2").Value = 1
$ws2.Range("X2").Value = 1
$ws2.Range("Y2").Value = 1
$ws2.Range("Z2").Value = 1
$ws2.Range("AA2").Value = 1
$ws2.Range("AB2").Value = 1
$ws2.Range("AC2").Value = 1
$ws2.Range("AD2").Value = 1
$ws2.Range("AE2").Value = 1
$ws2.Range("AF2").Value = 1
$ws2.Range("AG2").Value = 1
$ws2.Range("AH2").Value = 1
$ws2.Range("AI2").Value = 1
$ws2.Range("AJ2").Value = 1
$ws2.Range("AK2").Value = 1
$ws2.Range("AL2").Value = 1
$ws2.Range("AM2").Value = 1
$ws2.Range("AN2").Value = 1
$ws2.Range("AO2").Value = 1
$ws2.Range("AP2").Value = 1
$ws2.Range("AQ2").Value = 1
$ws2.Range("AR2").Value = 1
$ws2.Range("AS2").Value = 1
$ws2.Range("AT2").Value = 1
$ws2.Range("AU2").Value = 1
$ws2.Range("AV2").Value = 1

# Row 3 (Old Product B (SOC)): extend yearly cells into half-yearly steps
# across the same historical span (2018-2025.5)
$ws2.Range("L3").Value = 1
$ws2.Range("M3").Value = 1
$ws2.Range("N3").Value = 1
$ws2.Range("O3").Value = 1
$ws2.Range("P3").Value = 1
$ws2.Range("Q3").Value = 1
$ws2.Range("R3").Value = 1
$ws2.Range("S3").Value = 1

# ---------------------------------------------------------------
# View state: restore selections recorded in the saved workbook
# (MarketShare is the visible/active tab, selected last so it keeps
# tabSelected).
# ---------------------------------------------------------------
$ws1.Range("AZ17").Select()
$ws2.Range("Q3").Select()

